$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers - copy format from H1 (bold/border/centered style) then set values
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-17: I = 1, J = same as H
$hValues = @{
    2 = 5
    3 = 6
    4 = 3
    5 = 7
    6 = 5
    7 = 5
    8 = 2
    9 = 4
    10 = 5
    11 = 6
    12 = 5
    13 = 5
    14 = 2
    15 = 5
    16 = 5
    17 = 5
}

foreach ($r in 2..17) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hValues[$r]
}

# Row 18 is special: I18=5, J18=7
$ws.Cells.Item(18, 9).Value = 5
$ws.Cells.Item(18, 10).Value = 7
